# Add a new "2022-Q3" quarterly sheet to the workbook and record it in the
# "总计" (totals) summary sheet.
#
# Before: 总计, 2022-Q2, 2022-Q1, 2021-Q4, 2021-Q3, 2021-Q2
# After:  总计, 2022-Q3, 2022-Q2, 2022-Q1, 2021-Q4, 2021-Q3, 2021-Q2

$wb = $excel.ActiveWorkbook
$total = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1) Update the "总计" sheet: push the existing data rows (2..6) down by
#    one row (3..7), and write the new 2022-Q3 summary into row 2.
# ---------------------------------------------------------------------

# Final (post-edit) contents of rows 2..7, column B/C/D, in order.
$rows = @(
    @("2022-Q3", 3, 0.6),
    @("2022-Q2", 14, 2.59),
    @("2022-Q1", 12, 2.35),
    @("2021-Q4", 14, 5.08),
    @("2021-Q3", 18, 4.49),
    @("2021-Q2", 7, 1.13)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $total.Cells.Item($r, 1).Value = $i
    $total.Cells.Item($r, 2).Value = $rows[$i][0]
    $total.Cells.Item($r, 3).Value = $rows[$i][1]
    $total.Cells.Item($r, 4).Value = $rows[$i][2]
}

# Keep the A-column index style ("s=2", same as the header/borders used
# elsewhere in this sheet) consistent on the newly-added row 7.
$total.Range("A2").Copy()
$total.Range("A7").PasteSpecial(-4122)
$total.Cells.Item(7, 1).Value = 5

# ---------------------------------------------------------------------
# 2) Insert the brand-new "2022-Q3" worksheet right after "总计".
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($null, $total)
$q3.Name = "2022-Q3"

# Match the page margins used throughout the rest of the workbook.
$q3.PageSetup.LeftMargin = 54
$q3.PageSetup.RightMargin = 54
$q3.PageSetup.TopMargin = 72
$q3.PageSetup.BottomMargin = 72
$q3.PageSetup.HeaderMargin = 36
$q3.PageSetup.FooterMargin = 36

# Header row (B1:H1) - copy the header style ("s=2") from the 总计 sheet.
$total.Range("B1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)

$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

# Data rows 2..4. Columns B, C, D, E, F, G are stored as text in the source
# workbook (so that numeric-looking strings like "870009", "008135" or
# "6.90" keep their exact textual representation, incl. leading/trailing
# zeros); columns A and H are real numbers.
$q3.Range("B2:G4").NumberFormat = "@"

# Keep column-A index values styled like the rest of the workbook (s="2").
$total.Range("A2").Copy()
$q3.Range("A2:A4").PasteSpecial(-4122)

function Set-FundRow($row, $idx, $code, $name, $scale, $position, $ratio, $value, $rank) {
    $q3.Cells.Item($row, 1).Value = $idx
    $q3.Cells.Item($row, 2).Value = $code
    $q3.Cells.Item($row, 3).Value = $name
    $q3.Cells.Item($row, 4).Value = $scale
    $q3.Cells.Item($row, 5).Value = $position
    $q3.Cells.Item($row, 6).Value = $ratio
    $q3.Cells.Item($row, 7).Value = $value
    $q3.Cells.Item($row, 8).Value = $rank
}

Set-FundRow 2 0 "870009" "广发资管平衡精选一年持有混合A" "7.47" "92.14" "6.90" "0.5154" 6
Set-FundRow 3 1 "872019" "广发资管平衡精选一年持有混合C" "1.09" "92.14" "6.90" "0.0752" 6
Set-FundRow 4 2 "008135" "华宸未来价值先锋混合" "0.28" "87.82" "3.59" "0.0101" 10

$q3.Range("A1").Select()

# Restore the original active sheet / selection: the last sheet ("2021-Q2")
# was the active tab before this edit (Worksheets.Add() activates the
# freshly-inserted sheet as a side effect, so move the selection back).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Select()
$lastSheet.Range("A1").Select()
